$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the canonical URL and the publication Date ---
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B2").Value = "https://molic-avc.gabriellesantosleandro.com/StructureDefinition/molicavc-race-extension"
$wsMetadata.Range("B8").Value = "2023-08-16T00:27:03-03:00"

# --- Elements sheet: the same StructureDefinition URL is repeated (fixed value
#     for Extension.url), and the ValueSet URL used by the binding also moved
#     to the new subdomain-based host ---
$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("R4").Value = "https://molic-avc.gabriellesantosleandro.com/StructureDefinition/molicavc-race-extension"
$wsElements.Range("Z6").Value = "https://molic-avc.gabriellesantosleandro.com/ValueSet/molicavc-race"

# --- Elements sheet: column Z (26) got narrower after the binding text changed ---
$wsElements.Columns.Item(26).ColumnWidth = 63
